# Update the "include in DreamTeam ?" status for a few players on both
# sheets, then leave the workbook with MI as the active/selected tab
# (RR was previously the selected tab) and update each sheet's last
# on-screen cell selection.

$wb = $excel.ActiveWorkbook

$rr = $wb.Worksheets.Item("RR")
$mi = $wb.Worksheets.Item("MI")

# --- RR sheet ("RR" = Rajasthan Royals) value changes ---
# Y Jaiswal (row 6): Yes -> No
$rr.Range("E6").Value = "No"
# P Krishna (row 8): Yes -> No
$rr.Range("E8").Value = "No"
# Y Chahal (row 12): Yes -> Mandatory
$rr.Range("E12").Value = "Mandatory"

# --- MI sheet ("MI" = Mumbai Indians) value changes ---
# Rohit Sharma (row 4): No -> Mandatory
$mi.Range("E4").Value = "Mandatory"
# J Bumrah (row 12): Yes -> Mandatory
$mi.Range("E12").Value = "Mandatory"

# --- Selection/active-tab bookkeeping to match the saved UI state ---
# RR is no longer the selected tab; its last selection moved to F10.
$rr.Activate()
$rr.Range("F10").Select()

# MI becomes the selected/active tab, with E8 as the last selection.
$mi.Activate()
$mi.Range("E8").Select()
